$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row (66) values: B66, C66, D66
$ws.Range("B66").Value = 0.4
$ws.Range("C66").Value = 0.3
$ws.Range("D66").Value = 0.5

# Add new row 67 with the next month's data (01-07-2021)
# Format A67 as text first so the date-like "01-07-2021" string is
# stored as a literal string (matching column A's other entries)
# instead of being auto-converted to a date serial number.
$ws.Range("A67").NumberFormat = "@"
$ws.Range("A67").Value = "01-07-2021"
$ws.Range("A67").Style = "Normal"

$ws.Range("B67").Value = 1.3
$ws.Range("C67").Value = 0.5
$ws.Range("D67").Value = 1.4
